$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 25078.25
$ws.Range("I11").Value = 25078.25
$ws.Range("K11").Value = 25078.25
$ws.Range("M11").Value = -24938.25
$ws.Range("H39").Value = 1229.6
$ws.Range("I39").Value = 1541.2
$ws.Range("K39").Value = 4623.6
$ws.Range("M39").Value = -4327.6
$ws.Range("H43").Value = 9276342
$ws.Range("I43").Value = 33666.668
$ws.Range("K43").Value = 33666.668
$ws.Range("M43").Value = -33597.668
$ws.Range("N44").ClearContents()
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H74").Value = 3500
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3500
$ws.Range("N74").Value = -5372
$ws.Range("M77").ClearContents()
$ws.Range("H77").Value = 3500
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 17500
$ws.Range("N77").Value = -26860
$ws.Range("H107").Value = 1408.9231
$ws.Range("I107").Value = 1275.6471
$ws.Range("J107").Value = 1660.6666
$ws.Range("K107").Value = 1275.6471
$ws.Range("L107").Value = 1660.6666
$ws.Range("M107").Value = 644.3529000000001
$ws.Range("N107").Value = -5500.6666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7656
$ws.Range("I2").Value = 501.22223
$ws.Range("J2").Value = 20534.6
$ws.Range("K2").Value = 501.22223
$ws.Range("L2").Value = 20534.6
$ws.Range("M2").Value = -388.22223
$ws.Range("N2").Value = -20760.6
$ws.Range("H63").Value = 1800
$ws.Range("I63").Value = 1500
$ws.Range("J63").Value = 2100
$ws.Range("K63").Value = 1500
$ws.Range("L63").Value = 2100
$ws.Range("M63").Value = -814
$ws.Range("N63").Value = -3472
$ws.Range("H66").Value = 1800
$ws.Range("I66").Value = 1500
$ws.Range("J66").Value = 2100
$ws.Range("K66").Value = 7500
$ws.Range("L66").Value = 10500
$ws.Range("M66").Value = -4068
$ws.Range("N66").Value = -17364
$ws.Range("H102").Value = 9262067
$ws.Range("I102").Value = 11113887
$ws.Range("K102").Value = 11113887
$ws.Range("M102").Value = -11112265
$ws.Range("H110").Value = 1140.238
$ws.Range("I110").Value = 908.8
$ws.Range("K110").Value = 908.8
$ws.Range("M110").Value = 1136.2
$ws.Range("H116").Value = 7656
$ws.Range("I116").Value = 501.22223
$ws.Range("J116").Value = 20534.6
$ws.Range("K116").Value = 501.22223
$ws.Range("L116").Value = 20534.6
$ws.Range("M116").Value = 1792.77777
$ws.Range("N116").Value = -25122.6
$ws.Range("H132").Value = 1962.2433
$ws.Range("I132").Value = 1660.7241
$ws.Range("K132").Value = 4982.1723
$ws.Range("M132").Value = -2452.1723

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7656
$ws.Range("I3").Value = 501.22223
$ws.Range("J3").Value = 20534.6
$ws.Range("K3").Value = 501.22223
$ws.Range("L3").Value = 20534.6
$ws.Range("M3").Value = -387.22223
$ws.Range("N3").Value = -20762.6
$ws.Range("H80").Value = 610.13635
$ws.Range("I80").Value = 380.54544
$ws.Range("K80").Value = 380.54544
$ws.Range("M80").Value = 617.45456
$ws.Range("H83").Value = 610.13635
$ws.Range("I83").Value = 380.54544
$ws.Range("K83").Value = 1902.7272
$ws.Range("M83").Value = 3089.2728
$ws.Range("H94").Value = 10000404
$ws.Range("I94").Value = 12500338
$ws.Range("K94").Value = 12500338
$ws.Range("M94").Value = -12499887
$ws.Range("H99").Value = 35715588
$ws.Range("I99").Value = 50001156
$ws.Range("K99").Value = 50001156
$ws.Range("M99").Value = -49999658

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 47619876
$ws.Range("I16").Value = 50000780
$ws.Range("K16").Value = 50000780
$ws.Range("M16").Value = -50000493
$ws.Range("H22").Value = 582.2222
$ws.Range("I22").Value = 373.33334
$ws.Range("K22").Value = 373.33334
$ws.Range("M22").Value = -23.33334000000002
$ws.Range("M50").ClearContents()
$ws.Range("H50").Value = 21500
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 21500
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 21500
$ws.Range("N50").Value = -22750
$ws.Range("H58").Value = 1102.8096
$ws.Range("I58").Value = 917.13336
$ws.Range("K58").Value = 917.13336
$ws.Range("M58").Value = -714.13336
$ws.Range("H105").Value = 752.9
$ws.Range("I105").Value = 725.44446
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 725.44446
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 1021.55554
$ws.Range("N105").Value = -4494
$ws.Range("H113").Value = 47619876
$ws.Range("I113").Value = 50000780
$ws.Range("K113").Value = 50000780
$ws.Range("M113").Value = -49998610
$ws.Range("H132").Value = 2509.6667
$ws.Range("I132").Value = 1922.3846
$ws.Range("K132").Value = 5767.1538
$ws.Range("M132").Value = -3237.1538
$ws.Range("H134").Value = 1090.7188
$ws.Range("I134").Value = 909.0357
$ws.Range("K134").Value = 2727.1071
$ws.Range("M134").Value = -192.1071000000002
$ws.Range("H136").Value = 1102.8096
$ws.Range("I136").Value = 917.13336
$ws.Range("K136").Value = 2751.40008
$ws.Range("M136").Value = -201.4000800000003

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 697.8182
$ws.Range("I92").Value = 726.8570999999999
$ws.Range("J92").Value = 647
$ws.Range("K92").Value = 2180.5713
$ws.Range("L92").Value = 1941
$ws.Range("M92").Value = -932.5712999999996
$ws.Range("N92").Value = -4437

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 802.8182
$ws.Range("I107").Value = 859.8125
$ws.Range("K107").Value = 859.8125
$ws.Range("M107").Value = 1060.1875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 120
$ws.Range("I9").Value = 120
$ws.Range("K9").Value = 120
$ws.Range("M9").Value = 104
$ws.Range("H46").Value = 7691.846
$ws.Range("I46").Value = 998.5
$ws.Range("J46").Value = 10666.667
$ws.Range("K46").Value = 998.5
$ws.Range("L46").Value = 10666.667
$ws.Range("M46").Value = -810.5
$ws.Range("N46").Value = -11042.667
$ws.Range("H61").Value = 1523.6923
$ws.Range("I61").Value = 1257
$ws.Range("J61").Value = 1834.8334
$ws.Range("K61").Value = 1257
$ws.Range("L61").Value = 1834.8334
$ws.Range("M61").Value = -1055
$ws.Range("N61").Value = -2238.8334
$ws.Range("H113").Value = 1523.6923
$ws.Range("I113").Value = 1257
$ws.Range("J113").Value = 1834.8334
$ws.Range("K113").Value = 1257
$ws.Range("L113").Value = 1834.8334
$ws.Range("M113").Value = 913
$ws.Range("N113").Value = -6174.8334
$ws.Range("N121").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M81").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
